# Yardımcı Ekipmanlar - Paralel Sesli Anons Model-2 olarak düzeltildi.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Row 3 describes "Paralel Sesli Anons" (Kabin - Kapı Üstü) which had the
# wrong product code / model. Correct the product code suffix from -01 to
# -02 and update the Model column accordingly.
$ws.Range("B3").Value = "EQ-VOA-00-000-CLI-P1B0-02"
$ws.Range("J3").Value = "Model-02"

# Leave the final selection where the author left it after editing.
$ws.Range("K11").Select()
